# Update the Metadata sheet (sheet 1) and Elements sheet (sheet 2) to reflect
# the new StructureDefinition publication metadata (version bump, new date,
# publisher/jurisdiction instead of the old "Contact" row, and the
# root-element Short/Definition text for the Elements table).

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# --- Metadata sheet -------------------------------------------------------

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank; now populated.
$meta.Range("B9").Value = "Alvearie Team"

# The old sheet had a duplicated "Contact" / "No display for ContactDetail"
# row (rows 10 and 11 were identical). Remove one of those rows entirely
# (which shifts everything below up by one, also removing the now-blank
# trailing row 21), then turn the remaining row into the new
# "Jurisdiction" / "United States of America" pair.
$meta.Rows.Item(10).Delete()

$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# --- Elements sheet --------------------------------------------------------

# Root element's Short/Definition text now describes this specific
# extension instead of the generic Extension boilerplate.
$elements.Range("K2").Value = "Assigned Practitioner"
$elements.Range("L2").Value = "The practitioner assigned to the patient"

"Done applying StructureDefinition metadata updates."
